$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 27 (2026-02) stats
$ws.Range("B27").Value = 6555
$ws.Range("C27").Value = 1018
$ws.Range("D27").Value = 6115027
$ws.Range("E27").Value = 932.8797864225781
$ws.Range("F27").Value = 10.16806722689076
$ws.Range("G27").Value = 7.4973600844773
$ws.Range("H27").Value = 25.42340724769752
